# Updating filtered feeds from workflow
# Appends a new feed-item row (row 21) to the "Filtered Feeds" worksheet,
# mirroring the layout/formatting used by the existing rows (link / keywords / title).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newLink     = "https://www.360dx.com/immunoassays/randox-laboratories-gains-fda-de-novo-clearance-hemophilia-cdx"
$newKeyword  = "CDx"
$newTitle    = "Randox Laboratories Gains FDA De Novo Clearance for Hemophilia CDx"

# Column A: link, formatted/styled like the other link cells (Hyperlink style)
$ws.Cells.Item(21, 1).Value = $newLink
$ws.Hyperlinks.Add($ws.Cells.Item(21, 1), $newLink)
$ws.Cells.Item(21, 1).Style = "Hyperlink"

# Column B: keyword
$ws.Cells.Item(21, 2).Value = $newKeyword

# Column C: title
$ws.Cells.Item(21, 3).Value = $newTitle

Write-Output "Added row 21 to Filtered Feeds"
